$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Replace two existing cue/target pairs with new mediator words ---
$ws.Range("A67").Value = "Surfen"
$ws.Range("B67").Value = "Strand"

$ws.Range("A60").Value = "Kastanien"
$ws.Range("B60").Value = "Rösten"

# --- Append a brand-new cue/target pair as row 80 ---
# Copy the formatting of the preceding data row (rows 52:79 use a tinted font style)
$ws.Range("A79:B79").Copy()
$ws.Range("A80:B80").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A80").Value = "Nudeln"
$ws.Range("B80").Value = "Pizza"

# --- Extend the conditional formatting (duplicate-check) ranges down to the new last row (80) ---
# The two rules scoped to column A only
$fcsA = $ws.Range("A2").FormatConditions
for ($i = 1; $i -le $fcsA.Count; $i++) {
    $fc = $fcsA.Item($i)
    if ($fc.AppliesTo.Address() -eq "`$A`$2:`$A`$79") {
        $fc.ModifyAppliesToRange($ws.Range("A2:A80"))
    }
}

# The column-B-only rule and the combined A:B rule swap/extend their applied ranges
$fcsB = $ws.Range("B2").FormatConditions
$bOnlyRule = $fcsB.Item(1)   # originally B2:B79
$sharedRule = $fcsB.Item(2)  # originally A2:B79

$bOnlyRule.ModifyAppliesToRange($ws.Range("A2:B80"))
$sharedRule.ModifyAppliesToRange($ws.Range("B2:B80"))
$bOnlyRule.Priority = 4
$sharedRule.Priority = 1

# --- Update the view: scroll so row 57 is near the top, select the new last cell ---
$excel.ActiveWindow.ScrollRow = 57
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A80").Select()
